$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '51.026.15'
$ws.Range('E2').Value2 = '  -2.04%  '
$ws.Range('D3').Value2 = '2.906.00'
$ws.Range('E3').Value2 = '  -2.29%  '
$ws.Range('E4').Value2 = '  -0.18%  '
$ws.Range('D5').Value2 = '''370.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  +4.59%  '
$ws.Range('D6').Value2 = '''101.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  -5.54%  '
$ws.Range('E7').Value2 = '  -3.69%  '
$ws.Range('D8').Value2 = '''1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -0.08%  '
$ws.Range('E9').Value2 = '  -4.65%  '
$ws.Range('D10').Value2 = '''36.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -3.87%  '
$ws.Range('E11').Value2 = '  +0.52%  '
$ws.Range('E12').Value2 = '  -2.52%  '
$ws.Range('D13').Value2 = '''18.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value2 = '  -4.82%  '
$ws.Range('D14').Value2 = '3.362.81'
$ws.Range('E14').Value2 = '  -2.48%  '
$ws.Range('D15').Value2 = '''7.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  -3.75%  '
$ws.Range('D16').Value2 = '2.912.33'
$ws.Range('E16').Value2 = '  -1.86%  '
$ws.Range('D17').Value2 = '''0.918'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  -8.85%  '
$ws.Range('D18').Value2 = '50.970.64'
$ws.Range('E18').Value2 = '  -2.19%  '
$ws.Range('E19').Value2 = '  -7.78%  '
$ws.Range('D20').Value2 = '''7.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  -4.02%  '
$ws.Range('D21').Value2 = '''12.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  -5.18%  '
$ws.Range('E22').Value2 = '  -3.35%  '
$ws.Range('D23').Value2 = '''68.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value2 = '  -2.20%  '
$ws.Range('D24').Value2 = '''258.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value2 = '  -1.94%  '
$ws.Range('D25').Value2 = '''2.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  -2.51%  '
$ws.Range('E26').Value2 = '  -0.03%  '
$ws.Range('D27').Value2 = '''0.166'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value2 = '  -6.47%  '
$ws.Range('D28').Value2 = '''25.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value2 = '  -4.75%  '
$ws.Range('D29').Value2 = '''7.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value2 = '  -6.77%  '
$ws.Range('E30').Value2 = '  -6.58%  '
$ws.Range('E31').Value2 = '  +1.71%  '
$ws.Range('D32').Value2 = '''9.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -4.45%  '
$ws.Range('E33').Value2 = '  -2.04%  '
$ws.Range('D34').Value2 = '''51.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value2 = '  +1.14%  '
$ws.Range('D35').Value2 = '''34.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value2 = '  -5.69%  '
$ws.Range('E36').Value2 = '  +0.45%  '
$ws.Range('D37').Value2 = '''0.0421'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  -5.33%  '
$ws.Range('E38').Value2 = '  -7.21%  '
$ws.Range('E39').Value2 = '  -4.90%  '
$ws.Range('D40').Value2 = '''2.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value2 = '  -4.65%  '
$ws.Range('D41').Value2 = '''1.83'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  -6.49%  '
$ws.Range('E42').Value2 = '  -3.47%  '
$ws.Range('B43').Value2 = 'EnergySwap'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value2 = '''21.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value2 = '  -2.82%  '
$ws.Range('B44').Value2 = 'Monero'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value2 = '''119.27'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -2.27%  '
$ws.Range('E45').Value2 = '  -1.27%  '
$ws.Range('D46').Value2 = '2.016.64'
$ws.Range('E46').Value2 = '  -4.69%  '
$ws.Range('E47').Value2 = '  -1.15%  '
$ws.Range('D48').Value2 = '''3.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value2 = '  -7.16%  '
$ws.Range('D49').Value2 = '3.197.20'
$ws.Range('E49').Value2 = '  -2.28%  '
$ws.Range('E50').Value2 = '  -1.57%  '
$ws.Range('D51').Value2 = '''0.0307'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value2 = '  -8.64%  '
